$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) sheet ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 5673
$ws.Range("G4").Value = 45
$ws.Range("G5").Value = "不可售"
$ws.Range("F9").Value = 4474
$ws.Range("F10").Value = 1802
$ws.Range("F12").Value = 121
$ws.Range("F13").Value = 3016
$ws.Range("F15").Value = 575
$ws.Range("F16").Value = 228
$ws.Range("F17").Value = 566
$ws.Range("F19").Value = 489
$ws.Range("F23").Value = 1270
$ws.Range("F25").Value = 1478
$ws.Range("F26").Value = 120
$ws.Range("F32").Value = 75
$ws.Range("F33").Value = 114
$ws.Range("F34").Value = 77
$ws.Range("F35").Value = 3251
$ws.Range("F36").Value = 728
$ws.Range("F38").Value = 192
$ws.Range("F40").Value = 1231

# --- 演出 (Performances) sheet ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 16

# --- 全部类型 (All types) sheet ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 5673
$ws.Range("G4").Value = 45
$ws.Range("G5").Value = "不可售"
$ws.Range("F9").Value = 4474
$ws.Range("F10").Value = 1802
$ws.Range("F12").Value = 121
$ws.Range("F13").Value = 3016
$ws.Range("F15").Value = 575
$ws.Range("F16").Value = 228
$ws.Range("F17").Value = 566
$ws.Range("F19").Value = 489
$ws.Range("F20").Value = 16
$ws.Range("F24").Value = 1270
$ws.Range("F26").Value = 1478
$ws.Range("F27").Value = 120
$ws.Range("F33").Value = 75
$ws.Range("F34").Value = 114
$ws.Range("F35").Value = 77
$ws.Range("F36").Value = 3251
$ws.Range("F38").Value = 728
$ws.Range("F40").Value = 192
$ws.Range("F42").Value = 1231
